# Add 'ongkosKirimBeli' column to the Produk sheet and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Produk")

# Insert a new column before F (between hargaLuarKota and modifiedDate)
# and give it the new header "ongkosKirimBeli".
$ws.Columns("F").Insert()
$ws.Range("F1").Value = "ongkosKirimBeli"
$ws.Range("F1").ColumnWidth = 14

# Make "Produk" the active sheet / tab, with F2 selected.
$ws.Activate()
$ws.Range("F2").Select()
